# Applies the cryptos-list price/volume refresh described in the commit diff.
# Target cells store plain text (inlineStr) values, e.g. "247.59" and "  -1.70%  ",
# not numbers -- so each write is forced to text with a leading apostrophe
# (classic Excel "treat as text" marker) and the cell style is reset to Normal
# afterwards so no stray "@"/quotePrefix formatting sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.183.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.67%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.246.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.85%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.08%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''247.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.70%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.15%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''77.66'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +5.33%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.03%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.629'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.25%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''42.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +7.80%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -1.31%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''7.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.49%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -2.69%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''2.582.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''14.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -2.89%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.865'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.80%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.236.75'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.89%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''42.033.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.78%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.0₃0985'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -2.13%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''6.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.57%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''72.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.84%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +4.02%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''232.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -2.41%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = '''11.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -2.00%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''3.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -6.50%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -5.07%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  +13.39%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''2.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +1.15%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''170.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.70%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''20.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -2.19%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''33.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +9.13%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.0831'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.37%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -4.52%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.10%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''4.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -1.69%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''4.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +2.79%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''14.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.40%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0302'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -2.54%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''5.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +0.37%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -6.55%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''113.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +13.01%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -5.51%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''60.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -1.43%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''8.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -5.13%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.100'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -3.02%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -0.22%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -2.81%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.06%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''2.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.21%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''4.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -13.30%  '
$ws.Range("E51").Style = "Normal"
